$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Suraj Sharma"
$ws.Range("D7").Copy()
$ws.Range("D4").PasteSpecial()
$ws.Range("E7").Copy()
$ws.Range("E4").PasteSpecial()
$ws.Range("F4").Value = "Network"

$excel.CutCopyMode = $false

$ws.Range("E21").Select()
